$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> column letter -> new value, applied as direct
# cell writes (source data has no formulas; values were refreshed from an
# external market-price feed and pasted as static numbers).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4127.7417
$ws.Range("I98").Value = 2634.6875
$ws.Range("J98").Value = 5720.3335
$ws.Range("K98").Value = 2634.6875
$ws.Range("L98").Value = 5720.3335
$ws.Range("M98").Value = -1136.6875
$ws.Range("N98").Value = -8716.333500000001
$ws.Range("H122").Value = 4127.7417
$ws.Range("I122").Value = 2634.6875
$ws.Range("J122").Value = 5720.3335
$ws.Range("K122").Value = 7904.0625
$ws.Range("L122").Value = 17161.0005
$ws.Range("M122").Value = -5454.0625
$ws.Range("N122").Value = -22061.0005
$ws.Range("H129").Value = 879.89
$ws.Range("I129").Value = 315
$ws.Range("J129").Value = 891.4184
$ws.Range("K129").Value = 945
$ws.Range("L129").Value = 2674.2552
$ws.Range("M129").Value = 4055
$ws.Range("N129").Value = -12674.2552
$ws.Range("H137").Value = 1538148.1
$ws.Range("I137").Value = 2071482.4
$ws.Range("J137").Value = 4812
$ws.Range("K137").Value = 6214447.199999999
$ws.Range("L137").Value = 14436
$ws.Range("M137").Value = -6211897.199999999
$ws.Range("N137").Value = -19536
$ws.Range("H138").Value = 2496.77
$ws.Range("I138").Value = 685.4666999999999
$ws.Range("J138").Value = 2816.4119
$ws.Range("K138").Value = 2056.4001
$ws.Range("L138").Value = 8449.235700000001
$ws.Range("M138").Value = 3083.5999
$ws.Range("N138").Value = -18729.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3761.1323
$ws.Range("I32").Value = 3175.5762
$ws.Range("J32").Value = 7599.778
$ws.Range("K32").Value = 3175.5762
$ws.Range("L32").Value = 7599.778
$ws.Range("M32").Value = -2888.5762
$ws.Range("N32").Value = -8173.778
$ws.Range("H61").Value = 2075.4707
$ws.Range("I61").Value = 1125.125
$ws.Range("J61").Value = 2920.2222
$ws.Range("K61").Value = 1125.125
$ws.Range("L61").Value = 2920.2222
$ws.Range("M61").Value = -913.125
$ws.Range("N61").Value = -3344.2222
$ws.Range("H74").Value = 5531.409
$ws.Range("I74").Value = 7708.0835
$ws.Range("J74").Value = 2919.4
$ws.Range("K74").Value = 7708.0835
$ws.Range("L74").Value = 2919.4
$ws.Range("M74").Value = -6834.0835
$ws.Range("N74").Value = -4667.4
$ws.Range("H77").Value = 5531.409
$ws.Range("I77").Value = 7708.0835
$ws.Range("J77").Value = 2919.4
$ws.Range("K77").Value = 38540.4175
$ws.Range("L77").Value = 14597
$ws.Range("M77").Value = -34172.4175
$ws.Range("N77").Value = -23333
$ws.Range("H102").Value = 1923.75
$ws.Range("I102").Value = 1627.1428
$ws.Range("K102").Value = 1627.1428
$ws.Range("M102").Value = -5.142800000000079
$ws.Range("H122").Value = 2169.2173
$ws.Range("I122").Value = 1566.2858
$ws.Range("K122").Value = 4698.857400000001
$ws.Range("M122").Value = -2248.857400000001
$ws.Range("H132").Value = 2521.5789
$ws.Range("I132").Value = 1074
$ws.Range("J132").Value = 5003.143
$ws.Range("K132").Value = 3222
$ws.Range("L132").Value = 15009.429
$ws.Range("M132").Value = -692
$ws.Range("N132").Value = -20069.429
$ws.Range("H136").Value = 2075.4707
$ws.Range("I136").Value = 1125.125
$ws.Range("J136").Value = 2920.2222
$ws.Range("K136").Value = 3375.375
$ws.Range("L136").Value = 8760.6666
$ws.Range("M136").Value = -825.375
$ws.Range("N136").Value = -13860.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 29250.715
$ws.Range("J82").Value = 35151.2
$ws.Range("L82").Value = 35151.2
$ws.Range("N82").Value = -35917.2
$ws.Range("H85").Value = 29250.715
$ws.Range("J85").Value = 35151.2
$ws.Range("L85").Value = 35151.2
$ws.Range("N85").Value = -37803.2
$ws.Range("H86").Value = 2067.6924
$ws.Range("I86").Value = 1808.8889
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 1808.8889
$ws.Range("L86").Value = 2650
$ws.Range("M86").Value = -685.8888999999999
$ws.Range("N86").Value = -4896
$ws.Range("H89").Value = 2067.6924
$ws.Range("I89").Value = 1808.8889
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 9044.4445
$ws.Range("L89").Value = 13250
$ws.Range("M89").Value = -3428.4445
$ws.Range("N89").Value = -24482
$ws.Range("H95").Value = 34736.844
$ws.Range("J95").Value = 34736.844
$ws.Range("L95").Value = 34736.844
$ws.Range("N95").Value = -40228.844
$ws.Range("H105").Value = 1663.4857
$ws.Range("I105").Value = 1643.8154
$ws.Range("J105").Value = 1919.2
$ws.Range("K105").Value = 1643.8154
$ws.Range("L105").Value = 1919.2
$ws.Range("M105").Value = 103.1846
$ws.Range("N105").Value = -5413.2
$ws.Range("H134").Value = 1807.027
$ws.Range("I134").Value = 1217.5
$ws.Range("J134").Value = 5580
$ws.Range("K134").Value = 3652.5
$ws.Range("L134").Value = 16740
$ws.Range("M134").Value = -1117.5
$ws.Range("N134").Value = -21810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781.037
$ws.Range("I31").Value = 953.7
$ws.Range("K31").Value = 953.7
$ws.Range("M31").Value = -658.7
$ws.Range("H34").Value = 2781.037
$ws.Range("I34").Value = 953.7
$ws.Range("K34").Value = 953.7
$ws.Range("M34").Value = -751.7
$ws.Range("H58").Value = 2686.5293
$ws.Range("I58").Value = 1633.037
$ws.Range("J58").Value = 6750
$ws.Range("K58").Value = 1633.037
$ws.Range("L58").Value = 6750
$ws.Range("M58").Value = -1430.037
$ws.Range("N58").Value = -7156
$ws.Range("H132").Value = 2727.8635
$ws.Range("I132").Value = 1594.6875
$ws.Range("J132").Value = 5749.6665
$ws.Range("K132").Value = 4784.0625
$ws.Range("L132").Value = 17248.9995
$ws.Range("M132").Value = -2254.0625
$ws.Range("N132").Value = -22308.9995
$ws.Range("H134").Value = 6196.826
$ws.Range("I134").Value = 7133.0625
$ws.Range("K134").Value = 21399.1875
$ws.Range("M134").Value = -18864.1875
$ws.Range("H136").Value = 2686.5293
$ws.Range("I136").Value = 1633.037
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 4899.111
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -2349.111
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27854.545
$ws.Range("I4").Value = 75075
$ws.Range("J4").Value = 871.4286
$ws.Range("K4").Value = 225225
$ws.Range("L4").Value = 2614.2858
$ws.Range("M4").Value = -225113
$ws.Range("N4").Value = -2838.2858
$ws.Range("H23").Value = 227.57143
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 257.16666
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 771.4999799999999
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1241.49998
$ws.Range("H38").Value = 166.66667
$ws.Range("I38").Value = 120
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 360
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = -13
$ws.Range("N38").Value = -1294
$ws.Range("H97").Value = 494
$ws.Range("I97").Value = 665
$ws.Range("J97").Value = 425.6
$ws.Range("K97").Value = 1995
$ws.Range("L97").Value = 1276.8
$ws.Range("M97").Value = -1499
$ws.Range("N97").Value = -2268.8
$ws.Range("H113").Value = 3677049.5
$ws.Range("I113").Value = 611.2
$ws.Range("J113").Value = 6579500.5
$ws.Range("K113").Value = 1833.6
$ws.Range("L113").Value = 19738501.5
$ws.Range("M113").Value = 336.3999999999999
$ws.Range("N113").Value = -19742841.5
$ws.Range("H122").Value = 3131.4707
$ws.Range("I122").Value = 1061.1111
$ws.Range("J122").Value = 3876.8
$ws.Range("K122").Value = 9549.999900000001
$ws.Range("L122").Value = 34891.2
$ws.Range("M122").Value = -7099.999900000001
$ws.Range("N122").Value = -39791.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3515.5454
$ws.Range("I102").Value = 2083.875
$ws.Range("J102").Value = 7333.3335
$ws.Range("K102").Value = 2083.875
$ws.Range("L102").Value = 7333.3335
$ws.Range("M102").Value = -461.875
$ws.Range("N102").Value = -10577.3335
$ws.Range("H122").Value = 4217.909
$ws.Range("I122").Value = 1342.4286
$ws.Range("J122").Value = 9250
$ws.Range("K122").Value = 4027.2858
$ws.Range("L122").Value = 27750
$ws.Range("M122").Value = -1577.2858
$ws.Range("N122").Value = -32650
$ws.Range("H126").Value = 3237.67
$ws.Range("I126").Value = 2904.831
$ws.Range("J126").Value = 4351.9565
$ws.Range("K126").Value = 8714.493
$ws.Range("L126").Value = 13055.8695
$ws.Range("M126").Value = -6244.493
$ws.Range("N126").Value = -17995.8695
$ws.Range("H132").Value = 3244.4614
$ws.Range("I132").Value = 1799.3125
$ws.Range("J132").Value = 5556.7
$ws.Range("K132").Value = 5397.9375
$ws.Range("L132").Value = 16670.1
$ws.Range("M132").Value = -2867.9375
$ws.Range("N132").Value = -21730.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 25633.908
$ws.Range("I56").Value = 23999.125
$ws.Range("J56").Value = 29993.334
$ws.Range("K56").Value = 23999.125
$ws.Range("L56").Value = 29993.334
$ws.Range("M56").Value = -23308.125
$ws.Range("N56").Value = -31375.334
$ws.Range("H100").Value = 3285.7144
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
$ws.Range("H132").Value = 4470.6787
$ws.Range("I132").Value = 928.41174
$ws.Range("J132").Value = 9945.091
$ws.Range("K132").Value = 2785.23522
$ws.Range("L132").Value = 29835.273
$ws.Range("M132").Value = -255.23522
$ws.Range("N132").Value = -34895.273
$ws.Range("H136").Value = 5130.7617
$ws.Range("I136").Value = 1860.6666
$ws.Range("K136").Value = 5581.9998
$ws.Range("M136").Value = -3031.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2957.3333
$ws.Range("I126").Value = 1552.6666
$ws.Range("J126").Value = 5766.6665
$ws.Range("K126").Value = 4657.9998
$ws.Range("L126").Value = 17299.9995
$ws.Range("M126").Value = -2187.9998
$ws.Range("N126").Value = -22239.9995
$ws.Range("H132").Value = 7754550.5
$ws.Range("I132").Value = 1716.3462
$ws.Range("K132").Value = 5149.0386
$ws.Range("M132").Value = -2619.0386
$ws.Range("H136").Value = 3570.7812
$ws.Range("I136").Value = 1765.0416
$ws.Range("J136").Value = 8988
$ws.Range("K136").Value = 5295.1248
$ws.Range("L136").Value = 26964
$ws.Range("M136").Value = -2745.1248
$ws.Range("N136").Value = -32064
